# London_tube_lines v2.xlsx edit:
#  - Rename several line-sheets to proper-case names
#  - Switch the active/selected sheet from "WC" to "HammersmithCity"
#  - Turn off iterative calculation (iterateCount=200 -> none)

$wb = $excel.ActiveWorkbook

# --- Rename sheets (case / label fixes) -----------------------------------
$renames = @{
    "bakerloo"         = "Bakerloo"
    "central"          = "Central"
    "circle"           = "Circle"
    "district"         = "District"
    "LondonOverground" = "Overground"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# --- Switch the active sheet from WC to HammersmithCity --------------------
# This moves tabSelected from the WC sheetView to HammersmithCity's sheetView
# and updates the workbook-level activeTab index.
$wb.Worksheets.Item("HammersmithCity").Activate()

# --- Disable iterative calculation -----------------------------------------
$excel.Iteration = $false
